# optimize dql test case logic
# Rename test-user fixture names pro003/pro027 -> txnbtree_pro003/txnbtree_pro027
# across the three places each appears (User_used column F, the "create user"
# Op_sql in column H, and the "show create user" Query_sql in column I), plus
# the cosmetic view/formatting fallout (column F width, wrapped-row heights,
# active selection) that Excel recomputes as a consequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Row 4 (txnbt_protocol_003 - create/show user with host)
$ws.Range("F4").Value = "txnbtree_pro003"
$ws.Range("H4").Value = "create user 'txnbtree_pro003'@'172.20.3.15' identified by 'abc123'"
$ws.Range("I4").Value = "show create user 'txnbtree_pro003'@'172.20.3.15'"

# Row 28 (txnbt_protocol_027 - create/show user with ssl)
$ws.Range("F28").Value = "txnbtree_pro027"
$ws.Range("H28").Value = "create user 'txnbtree_pro027' identified by 'abc123' require ssl"
$ws.Range("I28").Value = "show create user 'txnbtree_pro027'"

# Column F widened to fit the longer "txnbtree_pro0xx" user names.
$ws.Columns("F:F").ColumnWidth = 16.5

# Wrapped Op_sql cells reflow given the wider column F, growing rows 13 and 30.
$ws.Rows(13).RowHeight = 94.5
$ws.Rows(30).RowHeight = 351

# Selection moved to I2.
$ws.Range("I2").Select()
